# The deck ships two DrawingML themes:
#   theme1.xml -> linked from the (only) Slide Master  -> currently "Integral" / "Red Violet"
#   theme2.xml -> linked from the (only) Notes Master   -> currently "Office Theme" / "Office"
#
# The target edit swaps the two themes' colour palettes (the font scheme and
# format scheme are already identical between the two themes, so only the
# 12-slot colour scheme actually needs to move). We recolour the Slide
# Master's theme (theme1.xml, reachable from the object model) to the
# "Office" palette that theme2.xml carries, slot for slot, using
# ThemeColorScheme.Colors(i).RGB - the documented, granular way to edit a
# theme's colours from the PowerPoint object model.

$p = $ppt.ActivePresentation
$master = $p.Slides.Item(1).Master
$colorScheme = $master.Theme.ThemeColorScheme

# Index order (matches MsoThemeColorSchemeIndex / a:clrScheme child order):
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2  7 accent3
#   8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
# Values are the "Office" theme's sRGB colours (000000, FFFFFF, 44546A,
# E7E6E6, 5B9BD5, ED7D31, A5A5A5, FFC000, 4472C4, 70AD47, 0563C1, 954F72)
# encoded as Windows COLORREF ints (0x00BBGGRR) for the .RGB setter.
$officePalette = @(
    0,        # dk1      000000
    16777215, # lt1      FFFFFF
    6968388,  # dk2      44546A
    15132391, # lt2      E7E6E6
    13998939, # accent1  5B9BD5
    3243501,  # accent2  ED7D31
    10855845, # accent3  A5A5A5
    49407,    # accent4  FFC000
    12874308, # accent5  4472C4
    4697456,  # accent6  70AD47
    12673797, # hlink    0563C1
    7491477   # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = $officePalette[$i - 1]
}
